$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# --- Row 7: LSE_URL (base url) ---
$ws.Range("A7").Value = "LSE_URL"
$ws.Range("B7").Value = "https://www.londonstockexchange.com/live-markets/market-data-dashboard/price-explorer"
$ws.Range("C7").Value = "The base url of the LSE."
$ws.Range("C7").WrapText = $true

# --- Row 8: LSE_API_URL ---
$ws.Range("A8").Value = "LSE_API_URL"
$ws.Range("B8").Value = "https://api.londonstockexchange.com/api/v1/components/refresh"
$ws.Range("C8").Value = "The url for making post api calls."

# --- Row 9: LSE_SEARCH_URL ---
$ws.Range("A9").Value = "LSE_SEARCH_URL"
$ws.Range("B9").Value = "https://www.londonstockexchange.com/live-markets/market-data-dashboard/price-explorer?categories=EQUITY&subcategories=1&q=[STOCKCODE]&world=quotes"
$ws.Range("C9").Value = "The url for searching stock values by stock code. Replace [STOCKCODE] with the proper stock code and open the url."
$ws.Range("C9").WrapText = $true

# Widen column B to fit the long URL values that were just entered.
$ws.Columns.Item(2).ColumnWidth = 150.5

# Update the selected cell to B7, matching where the new data was entered.
$ws.Range("B7").Select() | Out-Null
